$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Passenger2")

# --- New header columns (added before the row-2 values change, to match
#     shared-string insertion order: Email / PhoneNumber are typed first) ---
$ws2.Range("C1").Value = "Email"
$ws2.Range("D1").Value = "PhoneNumber"

# give C1/D1 the same bordered look already used by A1 (reuse existing style)
$ws2.Range("A1").Copy()
$ws2.Range("C1:D1").PasteSpecial(-4122)

# --- Update existing passenger row ---
$ws2.Range("A2").Value = "Gangesh"
$ws2.Range("B2").Value = "Jha"

# --- New data cells in row 2 ---
# give C2/D2 the bordered look from A2 first, then layer the hyperlink font on C2
$ws2.Range("A2").Copy()
$ws2.Range("C2:D2").PasteSpecial(-4122)

$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:ujha777@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ujha777@gmail.com")

$ws2.Range("D2").Value = 1234567890

# --- Column widths for the two new columns ---
$ws2.Columns.Item(3).ColumnWidth = 22.6
$ws2.Columns.Item(4).ColumnWidth = 10.2

# --- View state: Passenger1 keeps its whole used range selected once it is
#     no longer the active tab; Passenger2 becomes the active tab / sheet,
#     with C6 selected ---
$ws1 = $wb.Worksheets.Item("Passenger1")
$ws1.Range("A1:B2").Select()

$ws2.Range("C6").Select()
$ws2.Activate()

$wb.Save()
